{"js": "// The \"Notice u/s 94 BNSS, 2023\" heading paragraph is being removed from\n// this template (it belongs on the Money Transfer letter, not the Put on\n// Hold letter that this document represents).\nconst body = context.document.body;\n\n// Locate the paragraph by its exact text.\nconst results = body.search(\"Notice u/s 94 BNSS, 2023\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Grab the whole paragraph that contains the match and delete it\n  // entirely (this also removes its paragraph mark / <w:p> element,\n  // not just the text inside it).\n  const para = results.items[0].paragraphs.getFirst();\n  para.delete();\n  await context.sync();\n}\n", "ps1": "# The \"Notice u/s 94 BNSS, 2023\" heading paragraph is being removed from\n# this template (it belongs on the Money Transfer letter, not the Put on\n# Hold letter that this document represents).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Notice u/s 94 BNSS, 2023\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    # Expand the found range to the whole paragraph so the paragraph mark\n    # is included too, then delete it outright (removes the <w:p> entirely\n    # instead of leaving an empty paragraph behind).\n    $range.Expand(4) | Out-Null  # wdParagraph = 4\n    $range.Delete()\n}\n"}
